$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -6.619099999999992
$ws.Range("D3").Value = -6.418999999999994
$ws.Range("D5").Value = -8.305899999999998
$ws.Range("E7").Value = 12.37
$ws.Range("A9").Value = -20.31589999999997
$ws.Range("E9").Value = 13.06820000000001
$ws.Range("D11").Value = -8.254299999999995
$ws.Range("D12").Value = -8.408600000000002
$ws.Range("A13").Value = -21.96790000000003
$ws.Range("A16").Value = -19.95739999999998
$ws.Range("A18").Value = -23.16770000000002
$ws.Range("A20").Value = -22.15560000000003
$ws.Range("D21").Value = -7.528699999999999
$ws.Range("E21").Value = 13.35870000000001
